# Update the "想去人数" (interest count) column F for the events whose
# counts changed in this refresh of the generated data.
#
# Mapping of event name (column C) -> new value for column F.
$updates = @{
    "蜀山·银泰百货高新店-2024漫趣地带嘉年华（免费）" = 278
    "合肥·第十五届次元之门动漫游戏博览会"             = 7703
    "合肥·首届AT次元时代动漫游戏嘉年华"               = 5609
    "合肥·Holic动漫游戏展"                            = 463
    "合肥·乐帮•崩原铁绝only同人首展"                  = 72
    "合肥·W·A第五人格同人only2.0"                     = 250
    "合肥·第九届环形宇宙动漫游戏嘉年华"               = 214
    "合肥·MAX特摄同人only2.0"                         = 54
}

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count

    for ($r = 1; $r -le $rowCount; $r++) {
        $nameCell = $ws.Cells.Item($r, 3)   # Column C = 名称
        $name = $nameCell.Value()

        if ($null -ne $name -and $updates.ContainsKey([string]$name)) {
            $ws.Cells.Item($r, 6).Value = $updates[[string]$name]   # Column F = 想去人数
        }
    }
}
